# FD05-EPIS-Informe ProyectoFinal.docx -- "Implementacion fd01 doc #1"
#
# 1) Remove the stray "_GoBack" bookmark that Word drops at the last
#    edit position (first paragraph, around the logo picture).
# 2) Clear the placeholder header text ("Logo de Mi Empresa" / "Logo de
#    mi Cliente") so the document no longer shows the sample header.

$d = $word.ActiveDocument

# --- 1) delete the "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2) clear the default header's placeholder content --------------
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)
if ($header.Exists) {
    $header.Range.Delete()
}
